$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H96").Value = 1131.0769
$ws.Range("I96").Value = 291.5
$ws.Range("J96").Value = 3929.6667
$ws.Range("K96").Value = 874.5
$ws.Range("L96").Value = 11789.0001
$ws.Range("M96").Value = 498.5
$ws.Range("N96").Value = -14535.0001

$ws.Range("H100").Value = 2193.84
$ws.Range("I100").Value = 2115.7273
$ws.Range("K100").Value = 2115.7273
$ws.Range("M100").Value = -1574.7273

$ws.Range("H121").Value = 779.5
$ws.Range("J121").Value = 779.5
$ws.Range("L121").Value = 2338.5
$ws.Range("N121").Value = -5832.5

$ws.Range("H140").Value = 89999
$ws.Range("J140").Value = 89999
$ws.Range("L140").Value = 89999
$ws.Range("N140").Value = -100359

$ws.Range("H141").Value = 2000
$ws.Range("I141").Value = 2000
$ws.Range("K141").Value = 6000
$ws.Range("M141").Value = -820

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 481.5
$ws.Range("I5").Value = 500
$ws.Range("J5").Value = 426
$ws.Range("K5").Value = 500
$ws.Range("L5").Value = 426
$ws.Range("M5").Value = -388
$ws.Range("N5").Value = -650

$ws.Range("H61").Value = 6064.4
$ws.Range("I61").Value = 5830.5
$ws.Range("K61").Value = 5830.5
$ws.Range("M61").Value = -5618.5

$ws.Range("H110").Value = 1594.4546
$ws.Range("I110").Value = 1631
$ws.Range("K110").Value = 1631
$ws.Range("M110").Value = 414

$ws.Range("H132").Value = 827.4167
$ws.Range("I132").Value = 827.4167
$ws.Range("K132").Value = 2482.2501
$ws.Range("M132").Value = 47.7498999999998

$ws.Range("H136").Value = 6064.4
$ws.Range("I136").Value = 5830.5
$ws.Range("K136").Value = 17491.5
$ws.Range("M136").Value = -14941.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 481.5
$ws.Range("I4").Value = 500
$ws.Range("J4").Value = 426
$ws.Range("K4").Value = 500
$ws.Range("L4").Value = 426
$ws.Range("M4").Value = -385
$ws.Range("N4").Value = -656

$ws.Range("H134").Value = 1889.5
$ws.Range("I134").Value = 1889.5
$ws.Range("K134").Value = 5668.5
$ws.Range("M134").Value = -3133.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 1984.2307
$ws.Range("J7").Value = 1527.2222
$ws.Range("L7").Value = 1527.2222
$ws.Range("N7").Value = -1753.2222

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H6").Value = 81.25
$ws.Range("I6").Value = 81.25
$ws.Range("K6").Value = 243.75
$ws.Range("M6").Value = -130.75

$ws.Range("H50").Value = 155.89473
$ws.Range("I50").Value = 150.66667
$ws.Range("J50").Value = 250
$ws.Range("K50").Value = 452.00001
$ws.Range("L50").Value = 750
$ws.Range("M50").Value = 28.99998999999997
$ws.Range("N50").Value = -1712

$ws.Range("H53").Value = 155.89473
$ws.Range("I53").Value = 150.66667
$ws.Range("J53").Value = 250
$ws.Range("K53").Value = 452.00001
$ws.Range("L53").Value = 750
$ws.Range("M53").Value = 28.99998999999997
$ws.Range("N53").Value = -1712

$ws.Range("H132").Value = 1924.2174
$ws.Range("I132").Value = 1215.2727
$ws.Range("J132").Value = 2574.0833
$ws.Range("K132").Value = 10937.4543
$ws.Range("L132").Value = 23166.7497
$ws.Range("M132").Value = -8407.454299999999
$ws.Range("N132").Value = -28226.7497

$ws.Range("H140").Value = 1750
$ws.Range("J140").Value = 4000
$ws.Range("L140").Value = 12000
$ws.Range("N140").Value = -22360

$ws.Range("H141").Value = 8585.799999999999
$ws.Range("I141").Value = 3232.25
$ws.Range("J141").Value = 30000
$ws.Range("K141").Value = 9696.75
$ws.Range("L141").Value = 90000
$ws.Range("M141").Value = -4516.75
$ws.Range("N141").Value = -100360

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 837.6667
$ws.Range("I80").Value = 607.5
$ws.Range("J80").Value = 1298
$ws.Range("K80").Value = 607.5
$ws.Range("L80").Value = 1298
$ws.Range("M80").Value = 390.5
$ws.Range("N80").Value = -3294

$ws.Range("H83").Value = 837.6667
$ws.Range("I83").Value = 607.5
$ws.Range("J83").Value = 1298
$ws.Range("K83").Value = 3037.5
$ws.Range("L83").Value = 6490
$ws.Range("M83").Value = 1954.5
$ws.Range("N83").Value = -16474

$ws.Range("H132").Value = 3412.5715
$ws.Range("I132").Value = 3198
$ws.Range("K132").Value = 9594
$ws.Range("M132").Value = -7064

$ws.Range("H138").Value = 104499.5
$ws.Range("J138").Value = 104499.5
$ws.Range("L138").Value = 104499.5
$ws.Range("N138").Value = -114779.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 818.13336
$ws.Range("I16").Value = 818.13336
$ws.Range("K16").Value = 818.13336
$ws.Range("M16").Value = -648.13336

$ws.Range("H22").Value = 2915.9167
$ws.Range("I22").Value = 1499
$ws.Range("J22").Value = 3388.2222
$ws.Range("K22").Value = 1499
$ws.Range("L22").Value = 3388.2222
$ws.Range("M22").Value = -1204
$ws.Range("N22").Value = -3978.2222

$ws.Range("H27").Value = 2915.9167
$ws.Range("I27").Value = 1499
$ws.Range("J27").Value = 3388.2222
$ws.Range("K27").Value = 1499
$ws.Range("L27").Value = 3388.2222
$ws.Range("M27").Value = -1392
$ws.Range("N27").Value = -3602.2222

$ws.Range("H40").Value = 7631
$ws.Range("I40").Value = 5292.5713
$ws.Range("K40").Value = 5292.5713
$ws.Range("M40").Value = -5156.5713

$ws.Range("H132").Value = 4137.8237
$ws.Range("J132").Value = 6401
$ws.Range("L132").Value = 19203
$ws.Range("N132").Value = -24263

$ws.Range("H136").Value = 4821.6665
$ws.Range("I136").Value = 4499.2856
$ws.Range("J136").Value = 5950
$ws.Range("K136").Value = 13497.8568
$ws.Range("L136").Value = 17850
$ws.Range("M136").Value = -10947.8568
$ws.Range("N136").Value = -22950

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H70").Value = 25000
$ws.Range("I70").Value = 25000
$ws.Range("K70").Value = 25000
$ws.Range("M70").Value = -24685

$ws.Range("H73").Value = 25000
$ws.Range("I73").Value = 25000
$ws.Range("K73").Value = 25000
$ws.Range("M73").Value = -23908

$ws.Range("H107").Value = 474.8
$ws.Range("I107").Value = 393
$ws.Range("J107").Value = 699.75
$ws.Range("K107").Value = 1179
$ws.Range("L107").Value = 2099.25
$ws.Range("M107").Value = 741
$ws.Range("N107").Value = -5939.25

$ws.Range("H132").Value = 5393
$ws.Range("I132").Value = 5351.3335
$ws.Range("K132").Value = 16054.0005
$ws.Range("M132").Value = -13524.0005

$ws.Range("H136").Value = 3839.5557
$ws.Range("I136").Value = 3508.2856
$ws.Range("K136").Value = 10524.8568
$ws.Range("M136").Value = -7974.856800000001
